$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay a text cell (column D prices often
# "look like" plain numbers, e.g. "311.54", which Excel's COM Value setter
# would otherwise silently coerce to a float). Briefly force text format,
# assign, then clear the format again so no residual style index is left
# on the cell (keeps cell "s" attribute identical to the original file).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "28.230.57"
$ws.Range("E2").Value = "  +0.73%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.871.91"
$ws.Range("E3").Value = "  +3.88%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.10%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "311.54"
$ws.Range("E5").Value = "  +0.47%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.03%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.5024"
$ws.Range("E7").Value = "  -1.05%  "

# Row 8 - Cardano
Set-TextValue $ws.Range("D8") "0.3924"
$ws.Range("E8").Value = "  +1.81%  "

# Row 9 - Dogecoin
Set-TextValue $ws.Range("D9") "0.09608"
$ws.Range("E9").Value = "  +6.91%  "

# Row 10 - Polygon
Set-TextValue $ws.Range("D10") "1.139"
$ws.Range("E10").Value = "  +3.99%  "

# Row 11 - OKB
Set-TextValue $ws.Range("D11") "41.01"
$ws.Range("E11").Value = "  +0.48%  "

# Row 12 - Polkadot
Set-TextValue $ws.Range("D12") "6.482"
$ws.Range("E12").Value = "  +1.64%  "

# Row 13 - Solana
Set-TextValue $ws.Range("D13") "20.99"
$ws.Range("E13").Value = "  +3.14%  "

# Row 14 - WrappedEther
Set-TextValue $ws.Range("D14") "1.877.18"
$ws.Range("E14").Value = "  +4.86%  "

# Row 15 - BinanceUSD
Set-TextValue $ws.Range("D15") "1.002"
$ws.Range("E15").Value = "  +0.00%  "

# Row 16 - Chainlink
Set-TextValue $ws.Range("D16") "7.409"
$ws.Range("E16").Value = "  +1.48%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +1.12%  "

# Row 18 - Litecoin
Set-TextValue $ws.Range("D18") "93.13"
$ws.Range("E18").Value = "  +1.00%  "

# Row 19 - TRON
Set-TextValue $ws.Range("D19") "0.06624"
$ws.Range("E19").Value = "  +0.74%  "

# Row 20 - Avalanche
Set-TextValue $ws.Range("D20") "17.48"
$ws.Range("E20").Value = "  +1.43%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.05%  "

# Row 22 - Uniswap
Set-TextValue $ws.Range("D22") "6.125"
$ws.Range("E22").Value = "  +1.93%  "

# Row 23 - WrappedBTC
Set-TextValue $ws.Range("D23") "28.284.20"
$ws.Range("E23").Value = "  +0.85%  "

# Row 24 - Cosmos
$ws.Range("E24").Value = "  +2.53%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +2.74%  "

# Row 26 - LidoDAOToken
Set-TextValue $ws.Range("D26") "2.543"
$ws.Range("E26").Value = "  +5.72%  "

# Row 27 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D27") "2.079.45"
$ws.Range("E27").Value = "  +3.69%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  +4.29%  "

# Row 29 - Monero
Set-TextValue $ws.Range("D29") "157.99"
$ws.Range("E29").Value = "  -0.18%  "

# Row 30 - BitcoinCash
$ws.Range("E30").Value = "  +0.07%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  +1.93%  "

# Row 32 - Stellar
Set-TextValue $ws.Range("D32") "0.1054"
$ws.Range("E32").Value = "  -3.01%  "

# Row 33 - Filecoin
Set-TextValue $ws.Range("D33") "5.630"
$ws.Range("E33").Value = "  +1.03%  "

# Row 34 - HuobiToken
Set-TextValue $ws.Range("D34") "3.626"
$ws.Range("E34").Value = "  -0.46%  "

# Row 35 - Hedera
Set-TextValue $ws.Range("D35") "0.06752"
$ws.Range("E35").Value = "  -2.45%  "

# Row 36 - FraxShare
Set-TextValue $ws.Range("D36") "9.476"
$ws.Range("E36").Value = "  +5.58%  "

# Row 37 - VeChain
Set-TextValue $ws.Range("D37") "0.02391"
$ws.Range("E37").Value = "  +2.47%  "

# Row 38 - Algorand
Set-TextValue $ws.Range("D38") "0.2174"
$ws.Range("E38").Value = "  +0.35%  "

# Row 39 - now TheSandbox (was Aptos)
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D39") "0.6376"
$ws.Range("E39").Value = "  +4.29%  "

# Row 40 - now Aptos (was TheSandbox)
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D40") "11.47"
$ws.Range("E40").Value = "  +0.64%  "

# Row 41 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D41") "4.972"
$ws.Range("E41").Value = "  -0.38%  "

# Row 42 - TrustWalletToken
Set-TextValue $ws.Range("D42") "1.175"
$ws.Range("E42").Value = "  +1.98%  "

# Row 43 - Frax
Set-TextValue $ws.Range("D43") "1.001"
$ws.Range("E43").Value = "  +0.12%  "

# Row 44 - EnergySwap
Set-TextValue $ws.Range("D44") "13.60"
$ws.Range("E44").Value = "  +2.85%  "

# Row 45 - Decentraland
$ws.Range("E45").Value = "  +2.79%  "

# Row 46 - PancakeSwap
Set-TextValue $ws.Range("D46") "3.663"
$ws.Range("E46").Value = "  -1.15%  "

# Row 47 - WEMIXTOKEN
Set-TextValue $ws.Range("D47") "1.265"
$ws.Range("E47").Value = "  -2.23%  "

# Row 48 - Quant
Set-TextValue $ws.Range("D48") "123.79"
$ws.Range("E48").Value = "  -0.69%  "

# Row 49 - NEARProtocol
Set-TextValue $ws.Range("D49") "1.986"
$ws.Range("E49").Value = "  +2.72%  "

# Row 50 - EOS
Set-TextValue $ws.Range("D50") "1.196"
$ws.Range("E50").Value = "  +1.05%  "

# Row 51 - Cronos
Set-TextValue $ws.Range("D51") "0.06841"
$ws.Range("E51").Value = "  +1.64%  "
